$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update stakeholder names (column A) ---
$ws.Range("A3").Value = "MidtTrafiks kunder"
$ws.Range("A4").Value = "MidtTrafik"
$ws.Range("A5").Value = "Datatilsynet"
$ws.Range("A6").Value = "Pengeinstituttet"

# --- Fix typo: Pengeinstitutet -> Pengeinstituttet (Behov column, row 6) ---
$ws.Range("C6").Value = "Pengeinstituttet har behov for at brugerens bestilling og betaling foregår sikkert"

# --- Fix typo: brugerene -> brugerne (Overvejelser column, row 3) ---
$ws.Range("E3").Value = "Da en større del af brugerne af sådan et system vil være ældre eller handicappede skal brugergrænsefladen være simpel og nem at finde rundt i."

# --- Fix typo: potentiele -> potentielle (Overvejelser column, row 4) ---
$ws.Range("E4").Value = "I tilfælde af at potentielle brugere stadig kontakter MidtTrafik telefonisk for at bestille Flexture, så vil det være logisk at MidtTrafiks side af systemet også kan bruges til at registrere kørsler."

# --- Fix typo: Intressent -> Interessent (Overvejelser column, row 5) ---
$ws.Range("E5").Value = "Da vi som udviklere ikke er godkendt til at benytte CPR registret, falder denne Interessent uden for betydning af denne analyse."

# --- Adjust row heights ---
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 60
